$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2  = 3
    3  = 0
    4  = 5
    5  = 6
    6  = 1
    7  = 2
    8  = 4
    9  = 2
    10 = 5
    11 = 2
    12 = 11
    13 = 9
    14 = 7
    15 = 6
    16 = 2
    17 = 6
    18 = 2
    19 = 4
    20 = 5
    21 = 9
    22 = 4
    23 = 3
    24 = 3
    25 = 1
    26 = 2
    27 = 1
    28 = 1
    29 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
